$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2, shifting existing rows 2-4 down to 3-5
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the "slug" versions of the header
# labels in row 1 (lowercase, hyphenated, accent-stripped).
$ws.Range("A2").Value = "edad-grupos-quinquenales"
$ws.Range("B2").Value = "n-parados"
$ws.Range("C2").Value = "aragon"
$ws.Range("D2").Value = "provincia-codigo"
$ws.Range("E2").Value = "provincia-nombre"
$ws.Range("F2").Value = "sexo"
$ws.Range("G2").Value = "mes-y-ano"
